# Update TOTAL RUNS (B) and WICKETS (C) columns for player standings rows 2-89
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$standings = @(
    @(2, 7, 0),
    @(3, 34, 0),
    @(4, 52, 0),
    @(5, 39, 0),
    @(6, 17, 0),
    @(7, 15, 0),
    @(8, 0, 2),
    @(9, 17, 1),
    @(10, 14, 2),
    @(11, 9, 2),
    @(12, 7, 3),
    @(13, 36, 0),
    @(14, 2, 0),
    @(15, 0, 0),
    @(16, 64, 0),
    @(17, 26, 0),
    @(18, 2, 0),
    @(19, 17, 3),
    @(20, 25, 2),
    @(21, 37, 1),
    @(22, 0, 2),
    @(23, 3, 2),
    @(24, 0, 0),
    @(25, 0, 0),
    @(26, 0, 0),
    @(27, 0, 0),
    @(28, 0, 0),
    @(29, 0, 0),
    @(30, 0, 0),
    @(31, 0, 0),
    @(32, 0, 0),
    @(33, 0, 0),
    @(34, 0, 0),
    @(35, 0, 0),
    @(36, 0, 0),
    @(37, 0, 0),
    @(38, 0, 0),
    @(39, 0, 0),
    @(40, 0, 0),
    @(41, 0, 0),
    @(42, 0, 0),
    @(43, 0, 0),
    @(44, 0, 0),
    @(45, 0, 0),
    @(46, 0, 0),
    @(47, 0, 0),
    @(48, 0, 0),
    @(49, 0, 0),
    @(50, 0, 0),
    @(51, 0, 0),
    @(52, 0, 0),
    @(53, 0, 0),
    @(54, 0, 0),
    @(55, 0, 0),
    @(56, 0, 0),
    @(57, 0, 0),
    @(58, 0, 0),
    @(59, 0, 0),
    @(60, 0, 0),
    @(61, 0, 0),
    @(62, 0, 0),
    @(63, 0, 0),
    @(64, 0, 0),
    @(65, 0, 0),
    @(66, 0, 0),
    @(67, 0, 0),
    @(68, 0, 0),
    @(69, 0, 0),
    @(70, 0, 0),
    @(71, 0, 0),
    @(72, 0, 0),
    @(73, 0, 0),
    @(74, 0, 0),
    @(75, 0, 0),
    @(76, 0, 0),
    @(77, 0, 0),
    @(78, 0, 0),
    @(79, 0, 0),
    @(80, 0, 0),
    @(81, 0, 0),
    @(82, 0, 0),
    @(83, 0, 0),
    @(84, 0, 0),
    @(85, 0, 0),
    @(86, 0, 0),
    @(87, 0, 0),
    @(88, 0, 0),
    @(89, 0, 0)
)

foreach ($row in $standings) {
    $rowNum = $row[0]
    $runs = $row[1]
    $wickets = $row[2]
    $ws.Cells.Item($rowNum, 2).Value2 = $runs
    $ws.Cells.Item($rowNum, 3).Value2 = $wickets
}

# Update the active selection shown when the workbook is opened
$ws.Range("E7").Select()
